# Live trading results update:
#  - Trade #129 (MarketMaking) closes out (early_exit).
#  - Two new trades open: #158 (momentum, UP) and #159 (HighProbConvergence, DOWN).
#  - Summary / Strategy Status roll-up numbers updated accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.48
$summary.Range("B4").Value = 0.6
$summary.Range("B6").Value = 128
$summary.Range("B7").Value = 61
$summary.Range("B9").Value = 47.66

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 99.56999999999999
$status.Range("D6").Value = 48
$status.Range("E6").Value = -0.24
$status.Range("F6").Value = -0.43
$status.Range("G6").Value = 47.92

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #129 (MarketMaking) closes - row 130
$allTrades.Cells.Item(130, 7).Value = 0.99              # G - Exit Price
$allTrades.Cells.Item(130, 8).Value = "CLOSED"          # H - Status
$allTrades.Cells.Item(130, 9).Value = 1.0204             # I - P&L %
$allTrades.Cells.Item(130, 10).Value = 0.01              # J - P&L $
$allTrades.Cells.Item(130, 11).Value = 99.56999999999999 # K - Capital After
$allTrades.Cells.Item(130, 12).Value = "early_exit"      # L - Exit Reason
$allTrades.Cells.Item(130, 13).Value = 0.17              # M - Duration (min)

# New trade #158 (momentum, UP) - row 159
$allTrades.Cells.Item(159, 1).Value = 158
$allTrades.Cells.Item(159, 2).NumberFormat = "@"
$allTrades.Cells.Item(159, 2).Value = "2026-02-18"
$allTrades.Cells.Item(159, 3).Value = "00:34:01"
$allTrades.Cells.Item(159, 4).Value = "momentum"
$allTrades.Cells.Item(159, 5).Value = "UP"
$allTrades.Cells.Item(159, 6).Value = 0.98
$allTrades.Cells.Item(159, 8).Value = "OPEN"
$allTrades.Cells.Item(159, 9).Value = 0
$allTrades.Cells.Item(159, 10).Value = 0
$allTrades.Cells.Item(159, 11).Value = 99.23374292899115
$allTrades.Cells.Item(159, 13).Value = 0
$allTrades.Cells.Item(159, 14).Value = 0
$allTrades.Cells.Item(159, 15).Value = 0
$allTrades.Cells.Item(159, 16).Value = 0.9
$allTrades.Cells.Item(159, 17).Value = "Upward momentum: 1.980% over 10 samples"

# New trade #159 (HighProbConvergence, DOWN) - row 160
$allTrades.Cells.Item(160, 1).Value = 159
$allTrades.Cells.Item(160, 2).NumberFormat = "@"
$allTrades.Cells.Item(160, 2).Value = "2026-02-18"
$allTrades.Cells.Item(160, 3).Value = "00:34:02"
$allTrades.Cells.Item(160, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(160, 5).Value = "DOWN"
$allTrades.Cells.Item(160, 6).Value = 0.02
$allTrades.Cells.Item(160, 8).Value = "OPEN"
$allTrades.Cells.Item(160, 9).Value = 0
$allTrades.Cells.Item(160, 10).Value = 0
$allTrades.Cells.Item(160, 11).Value = 100.4130057263667
$allTrades.Cells.Item(160, 13).Value = 0
$allTrades.Cells.Item(160, 14).Value = 0
$allTrades.Cells.Item(160, 15).Value = 0
$allTrades.Cells.Item(160, 16).Value = 0.95
$allTrades.Cells.Item(160, 17).Value = "Mean reversion DOWN: price 1.58% above mean (z=2.00)"

# ---------------------------------------------------------------------------
# momentum sheet - new trade #158 - row 41
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(41, 1).Value = 158
$momentum.Cells.Item(41, 2).NumberFormat = "@"
$momentum.Cells.Item(41, 2).Value = "2026-02-18"
$momentum.Cells.Item(41, 3).Value = "00:34:01"
$momentum.Cells.Item(41, 4).Value = "momentum"
$momentum.Cells.Item(41, 5).Value = "UP"
$momentum.Cells.Item(41, 6).Value = 0.98
$momentum.Cells.Item(41, 8).Value = "OPEN"
$momentum.Cells.Item(41, 9).Value = 0
$momentum.Cells.Item(41, 10).Value = 0
$momentum.Cells.Item(41, 11).Value = 99.23374292899115
$momentum.Cells.Item(41, 12).Value = 0
$momentum.Cells.Item(41, 13).Value = 0
$momentum.Cells.Item(41, 14).Value = 0.9
$momentum.Cells.Item(41, 15).Value = "Upward momentum: 1.980% over 10 samples"
$momentum.Cells.Item(41, 17).Value = 0

# ---------------------------------------------------------------------------
# HighProbConvergence sheet - new trade #159 - row 22
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(22, 1).Value = 159
$hpc.Cells.Item(22, 2).NumberFormat = "@"
$hpc.Cells.Item(22, 2).Value = "2026-02-18"
$hpc.Cells.Item(22, 3).Value = "00:34:02"
$hpc.Cells.Item(22, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(22, 5).Value = "DOWN"
$hpc.Cells.Item(22, 6).Value = 0.02
$hpc.Cells.Item(22, 8).Value = "OPEN"
$hpc.Cells.Item(22, 9).Value = 0
$hpc.Cells.Item(22, 10).Value = 0
$hpc.Cells.Item(22, 11).Value = 100.4130057263667
$hpc.Cells.Item(22, 12).Value = 0
$hpc.Cells.Item(22, 13).Value = 0
$hpc.Cells.Item(22, 14).Value = 0.95
$hpc.Cells.Item(22, 15).Value = "Mean reversion DOWN: price 1.58% above mean (z=2.00)"
$hpc.Cells.Item(22, 17).Value = 0

# ---------------------------------------------------------------------------
# MarketMaking sheet - trade #129 closes - row 50
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(50, 7).Value = 0.99               # G - Exit Price
$mm.Cells.Item(50, 8).Value = "CLOSED"           # H - Status
$mm.Cells.Item(50, 9).Value = 1.0204              # I - P&L %
$mm.Cells.Item(50, 10).Value = 0.01               # J - P&L $
$mm.Cells.Item(50, 11).Value = 99.56999999999999  # K - Capital After
$mm.Cells.Item(50, 16).Value = "early_exit"       # P - Exit Reason
$mm.Cells.Item(50, 17).Value = 0.17               # Q - Duration (min)
